$wb = $excel.ActiveWorkbook

# Rename the three "U*" sheets to "PL_U*"
$wb.Worksheets.Item("U1a").Name = "PL_U1a"
$wb.Worksheets.Item("U1b").Name = "PL_U1b"
$wb.Worksheets.Item("U2b").Name = "PL_U2b"

# Make "PL_U2b" the active sheet/tab (was "U1a" before)
$wb.Worksheets.Item("PL_U2b").Activate()
